# Apply the "adjust for CHP" update to the GDPbES workbook.
#
# The "crude oil" row (row 15) on the GDPbES sheet previously mirrored the
# "petroleum" row (row 11) via shared formulas (=B11, =C11, ...), all of
# which evaluated to 0. This edit hardcodes the Guaranteed Dispatch
# Fraction for crude oil to 1 (100%) for every forecast year (2015-2050,
# columns B:AK), replacing the formulas with literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GDPbES")

# Select the GDPbES sheet (matches the authored tabSelected/selection change)
$ws.Activate()
$ws.Range("A15").Select()

# Replace the formulas in row 15 (crude oil) with hardcoded value 1
$ws.Range("B15:AK15").Value = 1
